# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Periodo Mora" (column E) values, most recent period first, and the
# updated "Valor Mora" (column F) amounts for each data row (16-26).
$rows = @(16,17,18,19,20,21,22,23,24,25,26)
$periodos = @("2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303")
$valores  = @(180000,180000,180000,180000,180000,180000,180000,180000,220000,220000,220000)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Range("E$r").Value = $periodos[$i]
    $ws.Range("F$r").Value = $valores[$i]
}
